# Recompute profit-margin figures (columns H:N) on each profession sheet.
# Source values come from a scheduled crafting-cost/market-price refresh;
# this script just pokes the refreshed numbers into the existing cells.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Item ID 5487)
$ws.Range("H9").Value = 576.1818
$ws.Range("I9").Value = 753.5
$ws.Range("K9").Value = 753.5
$ws.Range("M9").Value = -584.5
# Row 11 (Item ID 5533)
$ws.Range("H11").Value = 52534.24
$ws.Range("I11").Value = 52534.24
$ws.Range("K11").Value = 52534.24
$ws.Range("M11").Value = -52394.24
# Row 28 (Item ID 27772)
$ws.Range("H28").Value = 53194.4
$ws.Range("I28").Value = 78195.53999999999
$ws.Range("J28").Value = 6763.7144
$ws.Range("K28").Value = 78195.53999999999
$ws.Range("L28").Value = 6763.7144
$ws.Range("M28").Value = -77710.53999999999
$ws.Range("N28").Value = -7733.7144
# Row 33 (Item ID 5512)
$ws.Range("H33").Value = 9644459
$ws.Range("I33").Value = 39624.74
$ws.Range("K33").Value = 39624.74
$ws.Range("M33").Value = -39395.74
# Row 40 (Item ID 5505)
$ws.Range("H40").Value = 2912.3333
$ws.Range("I40").Value = 2516.1667
$ws.Range("J40").Value = 3308.5
$ws.Range("K40").Value = 2516.1667
$ws.Range("L40").Value = 3308.5
$ws.Range("M40").Value = -2341.1667
$ws.Range("N40").Value = -3658.5
# Row 41 (Item ID 5478)
$ws.Range("H41").Value = 838446.5600000001
$ws.Range("I41").Value = 2500127.2
$ws.Range("K41").Value = 2500127.2
$ws.Range("M41").Value = -2499687.2
# Row 57 (Item ID 43247)
$ws.Range("H57").Value = 83000
$ws.Range("I57").Value = 74500
$ws.Range("K57").Value = 223500
$ws.Range("M57").Value = -223001
# Row 58 (Item ID 4606)
$ws.Range("H58").Value = 1717.4166
$ws.Range("I58").Value = 319.875
$ws.Range("J58").Value = 4512.5
$ws.Range("K58").Value = 959.625
$ws.Range("L58").Value = 13537.5
$ws.Range("M58").Value = -809.625
$ws.Range("N58").Value = -13837.5
# Row 64 (Item ID 5506)
$ws.Range("H64").Value = 4320
$ws.Range("J64").Value = 4244.4443
$ws.Range("L64").Value = 4244.4443
$ws.Range("N64").Value = -4740.4443
# Row 67 (Item ID 5506)
$ws.Range("H67").Value = 4320
$ws.Range("J67").Value = 4244.4443
$ws.Range("L67").Value = 4244.4443
$ws.Range("N67").Value = -5960.4443
# Row 80 (Item ID 12605)
$ws.Range("H80").Value = 2842308.5
$ws.Range("I80").Value = 4546194
$ws.Range("J80").Value = 2499
$ws.Range("K80").Value = 13638582
$ws.Range("L80").Value = 7497
$ws.Range("M80").Value = -13637584
$ws.Range("N80").Value = -9493
# Row 83 (Item ID 12605)
$ws.Range("H83").Value = 2842308.5
$ws.Range("I83").Value = 4546194
$ws.Range("J83").Value = 2499
$ws.Range("K83").Value = 40915746
$ws.Range("L83").Value = 22491
$ws.Range("M83").Value = -40910754
$ws.Range("N83").Value = -32475
# Row 86 (Item ID 12603)
$ws.Range("H86").Value = 9552135
$ws.Range("I86").Value = 4920.375
$ws.Range("J86").Value = 15427344
$ws.Range("K86").Value = 4920.375
$ws.Range("L86").Value = 15427344
$ws.Range("M86").Value = -3797.375
$ws.Range("N86").Value = -15429590
# Row 89 (Item ID 12603)
$ws.Range("H89").Value = 9552135
$ws.Range("I89").Value = 4920.375
$ws.Range("J89").Value = 15427344
$ws.Range("K89").Value = 24601.875
$ws.Range("L89").Value = 77136720
$ws.Range("M89").Value = -18985.875
$ws.Range("N89").Value = -77147952
# Row 93 (Item ID 18043)
$ws.Range("H93").Value = 70601
$ws.Range("J93").Value = 70601
$ws.Range("L93").Value = 70601
$ws.Range("M93").Value = -75593
# Row 107 (Item ID 27766)
$ws.Range("H107").Value = 1949
$ws.Range("I107").Value = 1811.25
$ws.Range("K107").Value = 1811.25
$ws.Range("M107").Value = 108.75
# Row 111 (Item ID 27768)
$ws.Range("H111").Value = 43255.25
$ws.Range("I111").Value = 35043
$ws.Range("J111").Value = 48182.6
$ws.Range("K111").Value = 105129
$ws.Range("L111").Value = 144547.8
$ws.Range("M111").Value = -102062
$ws.Range("N111").Value = -150681.8
# Row 113 (Item ID 27775)
$ws.Range("H113").Value = 100004390
$ws.Range("J113").Value = 4416.3335
$ws.Range("L113").Value = 4416.3335
$ws.Range("N113").Value = -10924.3335
# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 50004636
$ws.Range("J137").Value = 5945.5
$ws.Range("L137").Value = 17836.5
$ws.Range("N137").Value = -22936.5
# Row 138 (Item ID 44169)
$ws.Range("H138").Value = 3463.162
$ws.Range("J138").Value = 4533.625
$ws.Range("L138").Value = 13600.875
$ws.Range("N138").Value = -23880.875

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Item ID 44147)
$ws.Range("H32").Value = 198726.77
$ws.Range("I32").Value = 219899.52
$ws.Range("J32").Value = 3937.4
$ws.Range("K32").Value = 219899.52
$ws.Range("L32").Value = 3937.4
$ws.Range("M32").Value = -219612.52
$ws.Range("N32").Value = -4511.4
# Row 41 (Item ID 2501)
$ws.Range("H41").Value = 17000.076
$ws.Range("I41").Value = 13500
$ws.Range("K41").Value = 13500
$ws.Range("M41").Value = -13086
# Row 61 (Item ID 43999)
$ws.Range("H61").Value = 1373.8948
$ws.Range("I61").Value = 947.13336
$ws.Range("K61").Value = 947.13336
$ws.Range("M61").Value = -735.13336
# Row 102 (Item ID 19945)
$ws.Range("H102").Value = 4825.6
$ws.Range("I102").Value = 1800.2858
$ws.Range("J102").Value = 11884.667
$ws.Range("K102").Value = 1800.2858
$ws.Range("L102").Value = 11884.667
$ws.Range("M102").Value = -178.2858000000001
$ws.Range("N102").Value = -15128.667
# Row 110 (Item ID 27708)
$ws.Range("H110").Value = 250038400
$ws.Range("I110").Value = 333334530
$ws.Range("J110").Value = 150000
$ws.Range("K110").Value = 333334530
$ws.Range("L110").Value = 150000
$ws.Range("M110").Value = -333332485
$ws.Range("N110").Value = -154090
# Row 122 (Item ID 36168)
$ws.Range("H122").Value = 5823
$ws.Range("I122").Value = 2449.6667
$ws.Range("J122").Value = 8714.429
$ws.Range("K122").Value = 7349.000100000001
$ws.Range("L122").Value = 26143.287
$ws.Range("M122").Value = -4899.000100000001
$ws.Range("N122").Value = -31043.287
# Row 132 (Item ID 43997)
$ws.Range("H132").Value = 1723.5714
$ws.Range("I132").Value = 1723.5714
$ws.Range("K132").Value = 5170.7142
$ws.Range("M132").Value = -2640.7142
# Row 133 (Item ID 41857)
$ws.Range("H133").Value = 79999
$ws.Range("I133").Value = 79999
$ws.Range("K133").Value = 79999
$ws.Range("M133").Value = -77469
# Row 136 (Item ID 43999)
$ws.Range("H136").Value = 1373.8948
$ws.Range("I136").Value = 947.13336
$ws.Range("K136").Value = 2841.40008
$ws.Range("M136").Value = -291.4000800000003

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 17 (Item ID 2393)
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 50
$ws.Range("K17").Value = 50
$ws.Range("M17").Value = 122
# Row 64 (Item ID 14184)
$ws.Range("H64").Value = 11589.777
$ws.Range("J64").Value = 14755.143
$ws.Range("L64").Value = 14755.143
$ws.Range("N64").Value = -15205.143
# Row 67 (Item ID 14184)
$ws.Range("H67").Value = 11589.777
$ws.Range("J67").Value = 14755.143
$ws.Range("L67").Value = 14755.143
$ws.Range("N67").Value = -16315.143
# Row 86 (Item ID 12526)
$ws.Range("H86").Value = 38464320
$ws.Range("I86").Value = 62502144
$ws.Range("K86").Value = 62502144
$ws.Range("M86").Value = -62501021
# Row 89 (Item ID 12526)
$ws.Range("H89").Value = 38464320
$ws.Range("I89").Value = 62502144
$ws.Range("K89").Value = 312510720
$ws.Range("M89").Value = -312505104
# Row 99 (Item ID 19943)
$ws.Range("H99").Value = 1496.1904
$ws.Range("I99").Value = 1471
$ws.Range("K99").Value = 1471
$ws.Range("M99").Value = 27
# Row 105 (Item ID 19947)
$ws.Range("H105").Value = 1665.4147
$ws.Range("I105").Value = 1550.7407
$ws.Range("K105").Value = 1550.7407
$ws.Range("M105").Value = 196.2592999999999
# Row 107 (Item ID 27706)
$ws.Range("H107").Value = 33493602
$ws.Range("I107").Value = 311366.56
$ws.Range("K107").Value = 311366.56
$ws.Range("M107").Value = -309446.56
# Row 134 (Item ID 43998)
$ws.Range("H134").Value = 2918.5
$ws.Range("I134").Value = 2335.5334
$ws.Range("K134").Value = 7006.600199999999
$ws.Range("M134").Value = -4471.600199999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Item ID 5367)
$ws.Range("H22").Value = 1309.6842
$ws.Range("I22").Value = 530
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 530
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -180
$ws.Range("N22").Value = -3699
# Row 31 (Item ID 44023)
$ws.Range("H31").Value = 3231.6206
$ws.Range("J31").Value = 3443.65
$ws.Range("L31").Value = 3443.65
$ws.Range("N31").Value = -4033.65
# Row 34 (Item ID 44023)
$ws.Range("H34").Value = 3231.6206
$ws.Range("J34").Value = 3443.65
$ws.Range("L34").Value = 3443.65
$ws.Range("N34").Value = -3847.65
# Row 58 (Item ID 44021)
$ws.Range("H58").Value = 2159.5715
$ws.Range("I58").Value = 1441.3334
$ws.Range("J58").Value = 2499.7896
$ws.Range("K58").Value = 1441.3334
$ws.Range("L58").Value = 2499.7896
$ws.Range("M58").Value = -1238.3334
$ws.Range("N58").Value = -2905.7896
# Row 99 (Item ID 36198)
$ws.Range("H99").Value = 2420.2
$ws.Range("I99").Value = 2368.6667
$ws.Range("J99").Value = 2454.5557
$ws.Range("K99").Value = 2368.6667
$ws.Range("L99").Value = 2454.5557
$ws.Range("M99").Value = -870.6667000000002
$ws.Range("N99").Value = -5450.5557
# Row 126 (Item ID 36198)
$ws.Range("H126").Value = 2420.2
$ws.Range("I126").Value = 2368.6667
$ws.Range("J126").Value = 2454.5557
$ws.Range("K126").Value = 7106.000100000001
$ws.Range("L126").Value = 7363.6671
$ws.Range("M126").Value = -4636.000100000001
$ws.Range("N126").Value = -12303.6671
# Row 132 (Item ID 44019)
$ws.Range("H132").Value = 2168.8462
$ws.Range("I132").Value = 1618.7333
$ws.Range("J132").Value = 4002.5557
$ws.Range("K132").Value = 4856.199900000001
$ws.Range("L132").Value = 12007.6671
$ws.Range("M132").Value = -2326.199900000001
$ws.Range("N132").Value = -17067.6671
# Row 134 (Item ID 44020)
$ws.Range("H134").Value = 2876.3948
$ws.Range("J134").Value = 3493.5881
$ws.Range("L134").Value = 10480.7643
$ws.Range("N134").Value = -15550.7643
# Row 136 (Item ID 44021)
$ws.Range("H136").Value = 2159.5715
$ws.Range("I136").Value = 1441.3334
$ws.Range("J136").Value = 2499.7896
$ws.Range("K136").Value = 4324.0002
$ws.Range("L136").Value = 7499.3688
$ws.Range("M136").Value = -1774.0002
$ws.Range("N136").Value = -12599.3688

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Item ID 4847)
$ws.Range("H2").Value = 15625678
$ws.Range("I2").Value = 1143.4445
$ws.Range("K2").Value = 6860.667
$ws.Range("M2").Value = -6747.667
# Row 5 (Item ID 43974)
$ws.Range("H5").Value = 777.25
$ws.Range("I5").Value = 366.8
$ws.Range("K5").Value = 1100.4
$ws.Range("M5").Value = -988.4000000000001
# Row 32 (Item ID 4731)
$ws.Range("H32").Value = 2055.8
$ws.Range("I32").Value = 5002
$ws.Range("J32").Value = 1319.25
$ws.Range("K32").Value = 15006
$ws.Range("L32").Value = 3957.75
$ws.Range("M32").Value = -14723
$ws.Range("N32").Value = -4523.75
# Row 55 (Item ID 4733)
$ws.Range("H55").Value = 811.4286
$ws.Range("I55").Value = 570
$ws.Range("J55").Value = 1133.3334
$ws.Range("K55").Value = 1710
$ws.Range("L55").Value = 3400.0002
$ws.Range("M55").Value = -1533
$ws.Range("N55").Value = -3754.0002
# Row 61 (Item ID 4727)
$ws.Range("H61").Value = 377.77777
$ws.Range("I61").Value = 85.71429000000001
$ws.Range("K61").Value = 257.14287
$ws.Range("M61").Value = -42.14287000000002
# Row 81 (Item ID 12843)
$ws.Range("H81").Value = 1246.8
$ws.Range("I81").Value = 808.6667
$ws.Range("J81").Value = 1904
$ws.Range("K81").Value = 2426.0001
$ws.Range("L81").Value = 5712
$ws.Range("M81").Value = -1303.0001
$ws.Range("N81").Value = -7958
# Row 84 (Item ID 12843)
$ws.Range("H84").Value = 1246.8
$ws.Range("I84").Value = 808.6667
$ws.Range("J84").Value = 1904
$ws.Range("K84").Value = 7278.0003
$ws.Range("L84").Value = 17136
$ws.Range("M84").Value = -1662.0003
$ws.Range("N84").Value = -28368
# Row 97 (Item ID 19846)
$ws.Range("H97").Value = 1162.4166
$ws.Range("I97").Value = 1498.8334
$ws.Range("J97").Value = 826
$ws.Range("K97").Value = 4496.5002
$ws.Range("L97").Value = 2478
$ws.Range("M97").Value = -4000.5002
$ws.Range("N97").Value = -3470
# Row 107 (Item ID 27838)
$ws.Range("H107").Value = 895.65
$ws.Range("J107").Value = 628.3333
$ws.Range("L107").Value = 1884.9999
$ws.Range("N107").Value = -5724.9999
# Row 135 (Item ID 43974)
$ws.Range("H135").Value = 777.25
$ws.Range("I135").Value = 366.8
$ws.Range("K135").Value = 3301.2
$ws.Range("M135").Value = -766.2000000000003

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Item ID 36169)
$ws.Range("H102").Value = 2431.7
$ws.Range("I102").Value = 1467.55
$ws.Range("K102").Value = 1467.55
$ws.Range("M102").Value = 154.45
# Row 113 (Item ID 27710)
$ws.Range("H113").Value = 2851.3076
$ws.Range("I113").Value = 1665.4
$ws.Range("J113").Value = 3592.5
$ws.Range("K113").Value = 1665.4
$ws.Range("L113").Value = 3592.5
$ws.Range("M113").Value = 504.5999999999999
$ws.Range("N113").Value = -7932.5
# Row 122 (Item ID 36182)
$ws.Range("H122").Value = 2870.625
$ws.Range("I122").Value = 2215.6667
$ws.Range("J122").Value = 3712.7144
$ws.Range("K122").Value = 6647.000100000001
$ws.Range("L122").Value = 11138.1432
$ws.Range("M122").Value = -4197.000100000001
$ws.Range("N122").Value = -16038.1432
# Row 126 (Item ID 36184)
$ws.Range("H126").Value = 6014.825
$ws.Range("I126").Value = 7863.0454
$ws.Range("J126").Value = 3755.889
$ws.Range("K126").Value = 23589.1362
$ws.Range("L126").Value = 11267.667
$ws.Range("M126").Value = -21119.1362
$ws.Range("N126").Value = -16207.667
# Row 132 (Item ID 44008)
$ws.Range("H132").Value = 316444.4
$ws.Range("I132").Value = 479058.06
$ws.Range("J132").Value = 6000.1816
$ws.Range("K132").Value = 1437174.18
$ws.Range("L132").Value = 18000.5448
$ws.Range("M132").Value = -1434644.18
$ws.Range("N132").Value = -23060.5448

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Item ID 36249)
$ws.Range("H7").Value = 31253374
$ws.Range("I7").Value = 71430930
$ws.Range("K7").Value = 71430930
$ws.Range("M7").Value = -71430818
# Row 40 (Item ID 36248)
$ws.Range("H40").Value = 3791.5833
$ws.Range("I40").Value = 2999.8572
$ws.Range("K40").Value = 2999.8572
$ws.Range("M40").Value = -2863.8572
# Row 46 (Item ID 5282)
$ws.Range("H46").Value = 2264.9656
$ws.Range("I46").Value = 1781.3636
$ws.Range("J46").Value = 2560.5
$ws.Range("K46").Value = 1781.3636
$ws.Range("L46").Value = 2560.5
$ws.Range("M46").Value = -1593.3636
$ws.Range("N46").Value = -2936.5
# Row 55 (Item ID 5284)
$ws.Range("H55").Value = 408.2
$ws.Range("I55").Value = 451.46667
$ws.Range("K55").Value = 451.46667
$ws.Range("M55").Value = -278.46667
# Row 61 (Item ID 27740)
$ws.Range("H61").Value = 365294.1
$ws.Range("I61").Value = 408569.4
$ws.Range("J61").Value = 4666.6665
$ws.Range("K61").Value = 408569.4
$ws.Range("L61").Value = 4666.6665
$ws.Range("M61").Value = -408367.4
$ws.Range("N61").Value = -5070.6665
# Row 113 (Item ID 27740)
$ws.Range("H113").Value = 365294.1
$ws.Range("I113").Value = 408569.4
$ws.Range("J113").Value = 4666.6665
$ws.Range("K113").Value = 408569.4
$ws.Range("L113").Value = 4666.6665
$ws.Range("M113").Value = -406399.4
$ws.Range("N113").Value = -9006.666499999999
# Row 122 (Item ID 36247)
$ws.Range("H122").Value = 3699.8572
$ws.Range("I122").Value = 2919.8667
$ws.Range("K122").Value = 8759.6001
$ws.Range("M122").Value = -6309.6001
# Row 126 (Item ID 36249)
$ws.Range("H126").Value = 31253374
$ws.Range("I126").Value = 71430930
$ws.Range("K126").Value = 214292790
$ws.Range("M126").Value = -214290320
# Row 128 (Item ID 34582)
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 132 (Item ID 44058)
$ws.Range("H132").Value = 5440.6943
$ws.Range("I132").Value = 3613.5715
$ws.Range("J132").Value = 7998.6665
$ws.Range("K132").Value = 10840.7145
$ws.Range("L132").Value = 23995.9995
$ws.Range("M132").Value = -8310.7145
$ws.Range("N132").Value = -29055.9995
# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 7039
$ws.Range("J136").Value = 7527.067
$ws.Range("L136").Value = 22581.201
$ws.Range("N136").Value = -27681.201

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 61 (Item ID 2854)
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 113 (Item ID 27752)
$ws.Range("H113").Value = 2913.4167
$ws.Range("I113").Value = 1795.6666
$ws.Range("K113").Value = 5386.9998
$ws.Range("M113").Value = -3216.9998
# Row 122 (Item ID 36208)
$ws.Range("H122").Value = 2390.5789
$ws.Range("J122").Value = 3918
$ws.Range("L122").Value = 11754
$ws.Range("N122").Value = -16654
# Row 126 (Item ID 36210)
$ws.Range("H126").Value = 2164.6667
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
